$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Finalized Version ready for Publishing:
# refresh the LastLoginDate timestamps for the three user rows.
$ws.Range("J2").Value = 45973.3999544907
$ws.Range("J3").Value = 45973.411712419
$ws.Range("J4").Value = 45973.3890095949

# Column J (LastLoginDate) was narrowed to match the other date columns.
$ws.Columns.Item(10).ColumnWidth = 14.5
